$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New titration samples run after replacing titrant/acid (commit: "ran three
# blue tank samples" + batch re-checks) -- six new rows of CRM accuracy data
# appended below the existing A1:F23 table, all against the CRM opened
# 20210314 (shared string index 9 reused automatically by the engine).

$dates  = 20210328,20210328,20210328,20210328,20210328,20210328
$crm    = 2766.9919084923999,2226.22004515256,2224.4256267383998,2224.2019706982401,2225.54963650217,2228.2924875591002
$batch  = 2225.4699999999998,2226.4699999999998,2227.4699999999998,2228.4699999999998,2229.4699999999998,2230.4699999999998
$batchNum = 180,180,180,180,180,180
$notes  = "CRM opened 20210314","CRM opened 20210314","CRM opened 20210314","CRM opened 20210314","CRM opened 20210314","CRM opened 20210314"

for ($i = 0; $i -lt 6; $i++) {
    $r = 24 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $crm[$i]
    $ws.Cells.Item($r, 3).Value = $batch[$i]
    $ws.Cells.Item($r, 5).Value = $batchNum[$i]
    $ws.Cells.Item($r, 6).Value = $notes[$i]
}

# Column D: "% off" -- fill as one shared formula across the new block,
# matching how the previous blocks (D4:D11, D12:D19) were filled.
$ws.Range("D24:D29").Formula = "=100*(B24-C24)/C24"

# Leave the view parked near the newly-added rows, like the author did.
$ws.Range("A19").Select() | Out-Null
$ws.Range("H28").Select() | Out-Null
